$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; B=$null; C=$null; D="28.420.92"; E="  -1.09%  "},
    @{Row=3; B=$null; C=$null; D="1.875.82"; E="  -1.75%  "},
    @{Row=4; B=$null; C=$null; D="1.009"; E="  -2.12%  "},
    @{Row=5; B=$null; C=$null; D="315.79"; E="  -1.55%  "},
    @{Row=6; B=$null; C=$null; D="1.010"; E="  -2.07%  "},
    @{Row=7; B=$null; C=$null; D="0.5107"; E="  -2.00%  "},
    @{Row=8; B=$null; C=$null; D="0.3947"; E="  -0.16%  "},
    @{Row=9; B=$null; C=$null; D="0.08401"; E="  +0.48%  "},
    @{Row=10; B=$null; C=$null; D="1.108"; E="  -2.62%  "},
    @{Row=11; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="6.255"; E="  -0.90%  "},
    @{Row=12; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.876.03"; E="  -2.17%  "},
    @{Row=13; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="20.46"; E="  -1.25%  "},
    @{Row=14; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="7.257"; E="  -1.03%  "},
    @{Row=15; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.012"; E="  -1.76%  "},
    @{Row=16; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.00001105"; E="  -0.76%  "},
    @{Row=17; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="91.05"; E="  -0.97%  "},
    @{Row=18; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.06714"; E="  -1.70%  "},
    @{Row=19; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="17.69"; E="  -1.75%  "},
    @{Row=20; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.010"; E="  -2.03%  "},
    @{Row=21; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.949"; E="  -2.57%  "},
    @{Row=22; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="28.510.86"; E="  -1.12%  "},
    @{Row=23; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="11.12"; E="  -1.54%  "},
    @{Row=24; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.262"; E="  -0.67%  "},
    @{Row=25; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.095.04"; E="  -1.89%  "},
    @{Row=26; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="160.96"; E="  -1.54%  "},
    @{Row=27; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="20.70"; E="  -1.66%  "},
    @{Row=28; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.377"; E="  -2.92%  "},
    @{Row=29; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="127.38"; E="  -0.32%  "},
    @{Row=30; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.1049"; E="  -1.88%  "},
    @{Row=31; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.049"; E="  -0.77%  "},
    @{Row=32; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.780"; E="  -3.48%  "},
    @{Row=33; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="3.611"; E="  -2.40%  "},
    @{Row=34; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02432"; E="  -1.73%  "},
    @{Row=35; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.06491"; E="  -2.83%  "},
    @{Row=36; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.2185"; E="  -1.88%  "},
    @{Row=37; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="8.916"; E="  -5.98%  "},
    @{Row=38; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.262"; E="  -0.49%  "},
    @{Row=39; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.192"; E="  -0.58%  "},
    @{Row=40; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="5.086"; E="  +1.04%  "},
    @{Row=41; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.6437"; E="  -2.22%  "},
    @{Row=42; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="11.16"; E="  -0.36%  "},
    @{Row=43; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="1.010"; E="  -2.15%  "},
    @{Row=44; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.6055"; E="  -1.95%  "},
    @{Row=45; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="13.08"; E="  -1.31%  "},
    @{Row=46; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="3.705"; E="  -1.62%  "},
    @{Row=47; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="2.013"; E="  -0.56%  "},
    @{Row=48; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="122.10"; E="  -0.99%  "},
    @{Row=49; B=$null; C=$null; D="1.193"; E="  -8.61%  "},
    @{Row=50; B="EOS"; C="https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D="1.206"; E="  -2.85%  "},
    @{Row=51; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.06825"; E="  -2.00%  "},
)

foreach ($item in $rows) {
    if ($item.B -ne $null) { $ws.Cells.Item($item.Row, 2).Value = $item.B }
    if ($item.C -ne $null) { $ws.Cells.Item($item.Row, 3).Value = $item.C }
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.E
}
